$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" (changed) date column C for rows 2-9 from 45170 to 45174
$ws.Range("C2:C9").Value = 45174
